# aggiornamento fino a 20/09/2021
# Appends rows 375-385 to Sheet1, continuing the daily series in columns A:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: date serial (days since 1899-12-30), col B, col C, col D
$data = @(
    @(44449, 3,  13, 131.8191036300953),
    @(44450, 18, 29, 294.0580004055972),
    @(44451, 3,  31, 314.337862502535),
    @(44452, 4,  30, 304.1979314540661),
    @(44453, 1,  31, 314.337862502535),
    @(44454, 1,  32, 324.4777935510039),
    @(44455, 2,  32, 324.4777935510039),
    @(44456, 2,  31, 314.337862502535),
    @(44457, 3,  16, 162.2388967755019),
    @(44458, 1,  14, 141.9590346785642),
    @(44459, 3,  13, 131.8191036300953)
)

$lastExistingRow = 374
$startRow = 375

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    # Copy the formatting (date number format, border, alignment) from the
    # last existing row's A cell so the new date cell keeps the same style.
    $ws.Range("A$lastExistingRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
